$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.379.87"
$ws.Range("E2").Value = "  +8.44%  "

# Row 3
$ws.Range("D3").Value = "3.622.93"
$ws.Range("E3").Value = "  +8.20%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.37%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.94%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "193.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.59%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.652"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.15%  "

# Row 8
$ws.Range("D8").Value = "3.617.87"
$ws.Range("E8").Value = "  +8.50%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.21%  "

# Row 10
$ws.Range("E10").Value = "  +6.93%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.666"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.11%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.27%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000297"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.55%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.43%  "

# Row 15
$ws.Range("D15").Value = "4.179.43"
$ws.Range("E15").Value = "  +6.73%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.36%  "

# Row 17
$ws.Range("D17").Value = "3.613.92"
$ws.Range("E17").Value = "  +6.64%  "

# Row 18
$ws.Range("D18").Value = "70.244.21"
$ws.Range("E18").Value = "  +8.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.59%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.121"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.55%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.94%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "498.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +18.88%  "

# Row 25
$ws.Range("E25").Value = "  +10.32%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.33%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.86%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.69%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.86%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +17.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.73%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "618.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.47%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.117"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.23%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0831"
$ws.Range("E36").Value = "  +13.83%  "

# Row 37
$ws.Range("E37").Value = "  +5.26%  "

# Row 38
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.30%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.88%  "

# Row 40
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.401"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.33%  "

# Row 42
$ws.Range("D42").Value = "3.335.82"
$ws.Range("E42").Value = "  +7.79%  "

# Row 43
$ws.Range("E43").Value = "  +12.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.43%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0446"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.16%  "

# Row 46
$ws.Range("E46").Value = "  +18.02%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.39%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.138"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.63%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.93%  "

# Row 50
$ws.Range("E50").Value = "  +6.29%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.63%  "
